$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) and Volume (E) data columns to text format so that
# numeric-looking strings (e.g. "0.9995", "36.00") are preserved exactly
# as text instead of being auto-converted to numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '25.710.87'
$ws.Range("E2").Value = '  -5.89%  '
$ws.Range("D3").Value = '1.807.19'
$ws.Range("E3").Value = '  -5.19%  '
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '274.96'
$ws.Range("E5").Value = '  -10.35%  '
$ws.Range("D6").Value = '0.9993'
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").Value = '0.5053'
$ws.Range("E7").Value = '  -6.80%  '
$ws.Range("D8").Value = '0.3513'
$ws.Range("E8").Value = '  -7.85%  '
$ws.Range("D9").Value = '44.32'
$ws.Range("E9").Value = '  -3.23%  '
$ws.Range("D10").Value = '0.06628'
$ws.Range("E10").Value = '  -9.33%  '
$ws.Range("D11").Value = '19.92'
$ws.Range("E11").Value = '  -9.80%  '
$ws.Range("D12").Value = '0.8327'
$ws.Range("E12").Value = '  -7.87%  '
$ws.Range("D13").Value = '0.07792'
$ws.Range("E13").Value = '  -4.89%  '
$ws.Range("D14").Value = '1.795.31'
$ws.Range("E14").Value = '  +32.21%  '
$ws.Range("D15").Value = '5.054'
$ws.Range("E15").Value = '  -5.69%  '
$ws.Range("D16").Value = '87.34'
$ws.Range("E16").Value = '  -8.81%  '
$ws.Range("D17").Value = '0.9986'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").Value = '13.89'
$ws.Range("E18").Value = '  -6.69%  '
$ws.Range("D19").Value = '0.9997'
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").Value = '0.000007975'
$ws.Range("E20").Value = '  -7.96%  '
$ws.Range("D21").Value = '25.772.26'
$ws.Range("E21").Value = '  -5.32%  '
$ws.Range("D22").Value = '4.715'
$ws.Range("E22").Value = '  -6.75%  '
$ws.Range("D23").Value = '9.947'
$ws.Range("E23").Value = '  -8.08%  '
$ws.Range("D24").Value = '6.052'
$ws.Range("D25").Value = '141.83'
$ws.Range("E25").Value = '  -4.49%  '
$ws.Range("D26").Value = '2.121'
$ws.Range("E26").Value = '  -8.35%  '
$ws.Range("E27").Value = '  -5.96%  '
$ws.Range("D28").Value = '16.91'
$ws.Range("E28").Value = '  -7.92%  '
$ws.Range("D29").Value = '108.44'
$ws.Range("E29").Value = '  -7.21%  '
$ws.Range("D30").Value = '4.324'
$ws.Range("E30").Value = '  -10.87%  '
$ws.Range("D31").Value = '4.189'
$ws.Range("E31").Value = '  -10.69%  '
$ws.Range("D32").Value = '0.08772'
$ws.Range("E32").Value = '  -4.66%  '
$ws.Range("D33").Value = '0.04789'
$ws.Range("E33").Value = '  -5.63%  '
$ws.Range("D34").Value = '0.7226'
$ws.Range("E34").Value = '  -13.05%  '
$ws.Range("D35").Value = '1.132'
$ws.Range("E35").Value = '  -7.82%  '
$ws.Range("D36").Value = '2.870'
$ws.Range("E36").Value = '  -4.67%  '
$ws.Range("D37").Value = '0.9991'
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("D38").Value = '3.032'
$ws.Range("E38").Value = '  -8.82%  '
$ws.Range("D39").Value = '0.01853'
$ws.Range("E39").Value = '  -7.31%  '
$ws.Range("D40").Value = '0.5168'
$ws.Range("E40").Value = '  -13.94%  '
$ws.Range("D41").Value = '2.277'
$ws.Range("E41").Value = '  -16.00%  '
$ws.Range("D42").Value = '0.9435'
$ws.Range("E42").Value = '  -12.46%  '
$ws.Range("D43").Value = '112.65'
$ws.Range("E43").Value = '  -3.00%  '
$ws.Range("D44").Value = '6.150'
$ws.Range("E44").Value = '  -7.86%  '
$ws.Range("D45").Value = '7.995'
$ws.Range("E45").Value = '  -13.90%  '
$ws.Range("D46").Value = '0.9989'
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("D47").Value = '0.4556'
$ws.Range("E47").Value = '  -11.84%  '
$ws.Range("D48").Value = '0.1376'
$ws.Range("E48").Value = '  -10.23%  '
$ws.Range("D49").Value = '9.248'
$ws.Range("E49").Value = '  -9.31%  '
$ws.Range("D50").Value = '36.00'
$ws.Range("E50").Value = '  -5.51%  '
$ws.Range("D51").Value = '1.490'
$ws.Range("E51").Value = '  -9.25%  '

# Restore the default cell style (the quick NumberFormat toggle above
# would otherwise leave a stray text-format style on these cells).
$dataRange.Style = "Normal"
